# edit.ps1 - applies the tooltip / rules-text updates described in the diff.
#
# Strategy: Word's InsertXML, in this runtime, inserts the runs carried by a
# <w:p> wrapper *before* the start of the target Range without deleting the
# Range's own content (it does not behave like a destructive "replace").
# So for every paragraph we need to change we:
#   1. Find() the full, unique plain-text of the paragraph to get its Range.
#   2. InsertXML a complete replacement <w:p>...</w:p> (same pPr + new runs)
#      at a zero-length Range collapsed to the found Range's Start.
#   3. Re-Find the original text (now shifted after our freshly inserted
#      paragraph) and Delete() it, which also removes the obsolete trailing
#      paragraph mark merge artifact, leaving one clean paragraph behind.

$d = $word.ActiveDocument

function Replace-Paragraph {
    param($SearchText, $NewParagraphXml)

    $rng = $d.Content
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $SearchText"
        return
    }
    $start = $rng.Start

    $ins = $d.Range($start, $start)
    $ins.InsertXML($NewParagraphXml)

    $searchRng = $d.Range($start, $d.Content.End)
    $found2 = $searchRng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found2) {
        Write-Host "RE-FIND FAILED: $SearchText"
        return
    }
    $searchRng.Delete()
}

$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Restoring: ... -> split trailing run into ": " + rest (rest gets en-US lang)
Replace-Paragraph `
    "Restoring: Defend with it to Heal and amount proportional to the  difference between attacking and defending card" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Restoring</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Defend with this card to gain health equal to the difference of values between the defending and defended cards,</w:t></w:r></w:p>')

# 2) Bounce Cards: pPr gets lang=en-US, trailing run split into 4 runs
Replace-Paragraph `
    "Bounce Cards: When Defending, this card does  Damage = Def number-atk num " `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Bounce Cards</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: When </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>d</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">efending, this card does  </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>damage  equal  to the difference of values between the defending and defended cards,</w:t></w:r></w:p>')

# 3) Burn X: ... -> split trailing run into ": " + rest (rest gets en-US lang)
Replace-Paragraph `
    "Burn X: Deal X damage to target enemy when played with" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Burn X</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Deal 1 damage for each burn modifier on the card,</w:t></w:r></w:p>')

# 4) Parry: ... -> split trailing run into 5 runs
Replace-Paragraph `
    "Parry: Turn around with this card to deal that cards number as damage" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Parry</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Reverse </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">with this card to deal that cards </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>value</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> as damage</w:t></w:r></w:p>')

# 5) Draw X: ... -> split trailing run into ": " + rest (rest gets en-US lang)
Replace-Paragraph `
    "Draw X: draws X cards when card is played" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Draw X</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Draws 1 card for each draw modifier on the card</w:t></w:r></w:p>')

# 6) Cripple X: ... -> split trailing run into 4 runs
Replace-Paragraph `
    "Cripple X: makes opponent discard X random cards on modified card played" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Cripple X</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>M</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">akes opponent discard </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>1 card for each cripple modifier on the card</w:t></w:r></w:p>')

# 7) Spikey X -> Spiky X, trailing run split in two
Replace-Paragraph `
    "Spikey X: make opponent take X damage when this card is defended" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:hanging="360" w:left="720"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Spiky X</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>When this card is defended deal 1 damage for each spiky modifier of the card to the defending player.</w:t></w:r></w:p>')

# 8) Grilled steak: "+50 max hp" -> bold "Doubles max health"
Replace-Paragraph `
    "Grilled steak: +50 max hp IMPLEMENTED" `
    ('<w:p ' + $W + '><w:pPr><w:pStyle w:val="normal1"/><w:ind w:hanging="0" w:left="0"/><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Grilled steak: </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>Doubles max health</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="en-US"/></w:rPr><w:t>IMPLEMENTED</w:t></w:r></w:p>')

Write-Host "DONE"
